# Update countries & provincias Spain
#
# The sheet ("Pais") is a snapshot of per-country COVID counters kept sorted
# by "Casos totales" (column B) descending. This refresh updates the
# counters that moved for a handful of countries; four of them (Etiopia,
# Madagascar, Belice, Islas Virgenes Britanicas) overtake their neighbour in
# the ranking and bubble up one row, shifting the rows in between down.
# Reproduce that precisely by writing the post-refresh content straight into
# the affected rows, then refresh the "last updated" banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Countries whose counters changed but keep their row position ----

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1817409
$ws.Range("C4").Value = 589
$ws.Range("E4").Value = 1176596
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 105575

# Paises Bajos (row 25)
$ws.Range("B25").Value = 46442
$ws.Range("C25").Value = 185
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 5956

# Suecia (row 28)
$ws.Range("B28").Value = 37542
$ws.Range("C28").Value = 429
$ws.Range("E28").Value = 28176

# Emiratos Arabes Unidos (row 30)
$ws.Range("B30").Value = 34557
$ws.Range("C30").Value = 661
$ws.Range("D30").Value = 17932
$ws.Range("E30").Value = 16361
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 264

# Portugal (row 31)
$ws.Range("B31").Value = 32500
$ws.Range("C31").Value = 297
$ws.Range("D31").Value = 19409
$ws.Range("E31").Value = 11681
$ws.Range("G31").Value = 14
$ws.Range("H31").Value = 1410

# Dinamarca (row 50)
$ws.Range("B50").Value = 11669
$ws.Range("C50").Value = 36
$ws.Range("D50").Value = 10362
$ws.Range("E50").Value = 733
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 574

# Kazajistan (row 55)
$ws.Range("E55").Value = 5598
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 40

# Croacia (row 89)
$ws.Range("D89").Value = 2072
$ws.Range("E89").Value = 71

# Republica de Macedonia (row 90)
$ws.Range("B90").Value = 2226
$ws.Range("C90").Value = 62
$ws.Range("D90").Value = 1552
$ws.Range("E90").Value = 541
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 133

# Sri Lanka (row 101)
$ws.Range("B101").Value = 1628
$ws.Range("C101").Value = 15
$ws.Range("E101").Value = 817

# Republica de Chipre (row 122)
$ws.Range("B122").Value = 943
$ws.Range("E122").Value = 136

# Mozambique (row 156)
$ws.Range("B156").Value = 254
$ws.Range("C156").Value = 10
$ws.Range("D156").Value = 91
$ws.Range("E156").Value = 161

# Gibraltar (row 163)
$ws.Range("B163").Value = 170
$ws.Range("C163").Value = 1
$ws.Range("E163").Value = 21

# ---- Re-ranked block: Etiopia jumps ahead of Albania/Hong Kong/Tunez/Letonia ----

$ws.Range("A111").Value = "Etiopia"
$ws.Range("B111").Value = 1172
$ws.Range("C111").Value = 109
$ws.Range("D111").Value = 209
$ws.Range("E111").Value = 952
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 3
$ws.Range("H111").Value = 11

$ws.Range("A112").Value = "Albania"
$ws.Range("B112").Value = 1136
$ws.Range("C112").Value = 14
$ws.Range("D112").Value = 872
$ws.Range("E112").Value = 231
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 33

$ws.Range("A113").Value = "Hong Kong"
$ws.Range("B113").Value = 1085
$ws.Range("C113").Value = 2
$ws.Range("D113").Value = 1037
$ws.Range("E113").Value = 44
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 4

$ws.Range("A114").Value = "Tunez"
$ws.Range("B114").Value = 1077
$ws.Range("C114").Value = 1
$ws.Range("D114").Value = 960
$ws.Range("E114").Value = 69
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 48

$ws.Range("A115").Value = "Letonia"
$ws.Range("B115").Value = 1066
$ws.Range("C115").Value = 1
$ws.Range("D115").Value = 745
$ws.Range("E115").Value = 297
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 24

# ---- Re-ranked block: Madagascar jumps ahead of Andorra/Nicaragua/Chad ----

$ws.Range("A127").Value = "Madagascar"
$ws.Range("B127").Value = 771
$ws.Range("C127").Value = 13
$ws.Range("D127").Value = 168
$ws.Range("E127").Value = 597
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 6

$ws.Range("A128").Value = "Principado de Andorra"
$ws.Range("B128").Value = 764
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 692
$ws.Range("E128").Value = 21
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 51

$ws.Range("A129").Value = "Nicaragua"
$ws.Range("B129").Value = 759
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 370
$ws.Range("E129").Value = 354
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 35

$ws.Range("A130").Value = "Republica del Chad"
$ws.Range("B130").Value = 759
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 470
$ws.Range("E130").Value = 224
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 65

# ---- Swap: Belice overtakes Santa Lucia (tie at 18 total cases) ----

$ws.Range("A200").Value = "Belice"
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 16
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 18
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

# ---- Swap: Islas Virgenes Britanicas overtakes Papua Nueva Guinea (tie at 8) ----

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# ---- Refresh the "last updated" banner ----
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 14:05"
